$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.455242931842804
$ws.Range("B1").Value = 0.2978412210941315
$ws.Range("C1").Value = 0.2321763932704926
$ws.Range("D1").Value = 0.2262540459632874
$ws.Range("E1").Value = 0.2455356568098068
